$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# 1. Insert a new worksheet "_set_SCENARIOS" as the first tab, before
#    the existing "_set_PRODUCTS" sheet.
# ------------------------------------------------------------------
$firstSheet = $wb.Worksheets.Item(1)
$scenarios = $wb.Worksheets.Add($firstSheet)
$scenarios.Name = "_set_SCENARIOS"

# Populate the new sheet with its header + data (matches the layout
# used by the other "_set_*" sheets in this workbook).
$scenarios.Range("A1").Value = "s_Names"
$scenarios.Range("A2").Value = "low energy"
$scenarios.Range("A3").Value = "middle energy"
$scenarios.Range("A4").Value = "high energy"

# Copy the header style (bold, bordered, centered) from an existing
# "_set_*" sheet's header cell so we reuse the existing style record
# instead of creating a new one.
$wb.Worksheets.Item("_set_PRODUCTS").Range("A1").Copy() | Out-Null
$scenarios.Range("A1").PasteSpecial(-4122) | Out-Null  # xlPasteFormats
$excel.CutCopyMode = $false

# ------------------------------------------------------------------
# 2. Restore/adjust the selection on each sheet and make the new
#    sheet the active tab.
# ------------------------------------------------------------------
$wb.Worksheets.Item("_set_PRODUCTS").Columns.Item(1).Select() | Out-Null
$wb.Worksheets.Item("_set_RESOURCES").Range("F26").Select() | Out-Null

$scenarios.Activate() | Out-Null
$scenarios.Range("B7").Select() | Out-Null
